$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.706.01"
$ws.Range("E2").Value = "'  +1.30%  "
$ws.Range("D3").Value = "'2.277.17"
$ws.Range("E3").Value = "'  -0.29%  "
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("D5").Value = "'310.06"
$ws.Range("E5").Value = "'  -2.89%  "
$ws.Range("D6").Value = "'103.50"
$ws.Range("E6").Value = "'  +2.54%  "
$ws.Range("E7").Value = "'  +0.13%  "
$ws.Range("E8").Value = "'  +0.38%  "
$ws.Range("D9").Value = "'0.597"
$ws.Range("E9").Value = "'  -0.91%  "
$ws.Range("D10").Value = "'38.66"
$ws.Range("E10").Value = "'  -0.79%  "
$ws.Range("D11").Value = "'0.0897"
$ws.Range("E11").Value = "'  -0.33%  "
$ws.Range("D12").Value = "'8.19"
$ws.Range("E12").Value = "'  -0.47%  "
$ws.Range("E13").Value = "'  +1.30%  "
$ws.Range("D14").Value = "'0.971"
$ws.Range("E14").Value = "'  +1.87%  "
$ws.Range("D15").Value = "'14.99"
$ws.Range("E15").Value = "'  -0.46%  "
$ws.Range("D16").Value = "'2.624.25"
$ws.Range("E16").Value = "'  -0.22%  "
$ws.Range("D17").Value = "'2.274.81"
$ws.Range("E17").Value = "'  -0.44%  "
$ws.Range("D18").Value = "'42.537.70"
$ws.Range("E18").Value = "'  +0.68%  "
$ws.Range("D19").Value = "'7.22"
$ws.Range("E19").Value = "'  -1.11%  "
$ws.Range("E20").Value = "'  -0.53%  "
$ws.Range("D21").Value = "'13.26"
$ws.Range("E21").Value = "'  +4.22%  "
$ws.Range("D22").Value = "'72.79"
$ws.Range("E22").Value = "'  +0.11%  "
$ws.Range("D23").Value = "'3.40"
$ws.Range("E23").Value = "'  -3.87%  "
$ws.Range("D24").Value = "'262.21"
$ws.Range("E24").Value = "'  -2.18%  "
$ws.Range("E25").Value = "'  -1.95%  "
$ws.Range("E26").Value = "'  +0.36%  "
$ws.Range("D27").Value = "'10.62"
$ws.Range("E27").Value = "'  -1.44%  "
$ws.Range("E28").Value = "'  -0.20%  "
$ws.Range("E29").Value = "'  +14.98%  "
$ws.Range("D30").Value = "'22.15"
$ws.Range("E30").Value = "'  -1.10%  "
$ws.Range("D31").Value = "'35.63"
$ws.Range("E31").Value = "'  -4.81%  "
$ws.Range("D32").Value = "'164.44"
$ws.Range("E32").Value = "'  +0.27%  "
$ws.Range("D33").Value = "'0.0854"
$ws.Range("E33").Value = "'  -1.53%  "
$ws.Range("E34").Value = "'  -1.95%  "
$ws.Range("D35").Value = "'2.55"
$ws.Range("E35").Value = "'  +3.61%  "
$ws.Range("E36").Value = "'  -2.91%  "
$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = "'  -1.46%  "
$ws.Range("E38").Value = "'  -1.82%  "
$ws.Range("D39").Value = "'3.70"
$ws.Range("E39").Value = "'  +0.54%  "
$ws.Range("D40").Value = "'2.71"
$ws.Range("E40").Value = "'  -1.84%  "
$ws.Range("E41").Value = "'  +1.89%  "
$ws.Range("D42").Value = "'97.78"
$ws.Range("E42").Value = "'  +6.80%  "
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "'  +0.19%  "
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = "'68.50"
$ws.Range("E44").Value = "'  +0.40%  "
$ws.Range("E45").Value = "'  +0.76%  "
$ws.Range("D46").Value = "'1.713.96"
$ws.Range("E46").Value = "'  +6.49%  "
$ws.Range("D47").Value = "'11.86"
$ws.Range("E47").Value = "'  -1.56%  "
$ws.Range("D48").Value = "'109.50"
$ws.Range("E48").Value = "'  -5.11%  "
$ws.Range("D49").Value = "'75.17"
$ws.Range("E49").Value = "'  -5.05%  "
$ws.Range("D50").Value = "'5.15"
$ws.Range("E50").Value = "'  -0.83%  "
$ws.Range("E51").Value = "'  -4.02%  "
